$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix typo in existing shared string value (A1 row's activity label on row2? Actually it's A2)
$ws.Range("A2").Value = "chart and presentation"

# Correct the Date End value for row 2 (raw date serial, avoids host fractional-day quirk)
$ws.Range("C2").Value = 43731

# Add new row of data
$ws.Range("A3").Value = "Helped with designing of forms "
$ws.Range("B3").Value = 43729
$ws.Range("B3").NumberFormat = "d-mmm-yy"
$ws.Range("C3").Value = 43732
$ws.Range("C3").NumberFormat = "mm-dd-yy"

$ws.Range("A3").Select()
